$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value2 = "69.311.17"
$c = $ws.Cells.Item(2, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.37%  "

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.945.16"
$c = $ws.Cells.Item(3, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.32%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.998"
$c = $ws.Cells.Item(4, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.25%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value2 = "493.37"
$c = $ws.Cells.Item(5, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.13%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value2 = "147.29"
$c = $ws.Cells.Item(6, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.12%  "

$c = $ws.Cells.Item(7, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.00%  "

$c = $ws.Cells.Item(8, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.00%  "

$c = $ws.Cells.Item(9, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.07%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.178"
$c = $ws.Cells.Item(10, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +5.08%  "

$c = $ws.Cells.Item(11, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.87%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value2 = "43.31"
$c = $ws.Cells.Item(12, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.58%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value2 = "10.45"
$c = $ws.Cells.Item(13, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -2.79%  "

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value2 = "4.568.59"
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.30%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.939.41"
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.06%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value2 = "14.31"
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -2.69%  "

$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.79%  "

$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +4.23%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value2 = "19.92"
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.46%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value2 = "69.309.25"
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.12%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value2 = "439.10"
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.99%  "

$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.18%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value2 = "14.57"
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -2.12%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value2 = "89.01"
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.36%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value2 = "12.05"
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +9.79%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.80"
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +4.09%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value2 = "11.13"
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -2.46%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value2 = "37.22"
$c = $ws.Cells.Item(28, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -4.17%  "

$c = $ws.Cells.Item(29, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -3.99%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value2 = "703.58"
$c = $ws.Cells.Item(30, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.34%  "

$c = $ws.Cells.Item(31, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.82%  "

$c = $ws.Cells.Item(32, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.12%  "

$c = $ws.Cells.Item(33, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.42%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.466"
$c = $ws.Cells.Item(34, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +16.91%  "

$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0₃0902"
$c = $ws.Cells.Item(35, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.34%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value2 = "61.83"
$c = $ws.Cells.Item(36, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.41%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value2 = "6.07"
$c = $ws.Cells.Item(37, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.58%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value2 = "40.87"
$c = $ws.Cells.Item(38, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -2.61%  "

$c = $ws.Cells.Item(39, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.33%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.999"
$c = $ws.Cells.Item(40, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.20%  "

$c = $ws.Cells.Item(41, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.12%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0490"
$c = $ws.Cells.Item(42, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.56%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value2 = "2.91"
$c = $ws.Cells.Item(43, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +0.37%  "

$c = $ws.Cells.Item(44, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -3.78%  "

$c = $ws.Cells.Item(45, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +2.16%  "

$c = $ws.Cells.Item(46, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +1.05%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.40"
$c = $ws.Cells.Item(47, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +7.65%  "

$c = $ws.Cells.Item(48, 5)
$c.NumberFormat = "@"
$c.Value2 = "  +5.89%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value2 = "3.38"
$c = $ws.Cells.Item(49, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -1.02%  "

$ws.Cells.Item(50, 2).Value2 = "Monero"
$ws.Cells.Item(50, 3).Value2 = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value2 = "144.54"
$c = $ws.Cells.Item(50, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -0.99%  "

$ws.Cells.Item(51, 2).Value2 = "BabyDogeCoin"
$ws.Cells.Item(51, 3).Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value2 = "0.0₆0341"
$c = $ws.Cells.Item(51, 5)
$c.NumberFormat = "@"
$c.Value2 = "  -3.47%  "
